$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 18: # / Port / Pin# / uC IO function / Signal name Reflow / Function Group / Function description / uC voltage range [V] / Noted
$ws.Cells.Item(18, 1).Value = 21
$ws.Cells.Item(18, 2).Value = "A"
$ws.Cells.Item(18, 3).Value = 4
$ws.Cells.Item(18, 4).Value = "GPIO"
$ws.Cells.Item(18, 5).Value = "~CS_TEMP"
$ws.Cells.Item(18, 6).Value = "TEMP"
$ws.Cells.Item(18, 7).Value = "Temperature SPI chipselect"
$ws.Cells.Item(18, 8).Value = "[0..3.3]"

$ws.Range("E19").Select()
